$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="27.616.10"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = '  +2.60%  '

$ws.Range("D3").Formula = '="1.854.31"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = '  +2.34%  '

$ws.Range("D5").Formula = '="321.62"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +3.26%  '

$ws.Range("D6").Formula = '="1.029"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +2.51%  '

$ws.Range("D7").Formula = '="0.4396"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  +2.64%  '

$ws.Range("D8").Formula = '="0.3790"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  +2.99%  '

$ws.Range("D9").Formula = '="0.07420"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  +2.27%  '

$ws.Range("D10").Formula = '="0.8789"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  +2.16%  '

$ws.Range("D11").Formula = '="21.64"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  +2.37%  '

$ws.Range("D12").Formula = '="1.864.18"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  -8.76%  '

$ws.Range("D13").Formula = '="5.528"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  +2.81%  '

$ws.Range("D14").Formula = '="6.703"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  +1.18%  '

$ws.Range("D15").Formula = '="0.07215"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  +4.77%  '

$ws.Range("D16").Formula = '="83.22"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  +3.19%  '

$ws.Range("D17").Formula = '="1.035"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  +3.06%  '

$ws.Range("D18").Formula = '="0.000009058"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  +1.55%  '

$ws.Range("E19").Value = '  +2.55%  '

$ws.Range("D20").Formula = '="15.44"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +1.81%  '

$ws.Range("D21").Formula = '="27.633.37"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  +2.52%  '

$ws.Range("D22").Formula = '="5.260"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +1.40%  '

$ws.Range("D23").Formula = '="11.39"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  +3.57%  '

$ws.Range("D24").Formula = '="158.04"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +2.75%  '

$ws.Range("D25").Formula = '="1.923"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +2.06%  '

$ws.Range("D26").Formula = '="18.76"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  +2.77%  '

$ws.Range("D27").Formula = '="1.979"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +5.05%  '

$ws.Range("D28").Formula = '="5.292"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  +1.60%  '

$ws.Range("D29").Formula = '="117.20"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +2.03%  '

$ws.Range("D30").Formula = '="0.09065"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  +1.57%  '

$ws.Range("D31").Formula = '="1.203"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  +3.69%  '

$ws.Range("D32").Formula = '="0.7629"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  +2.85%  '

$ws.Range("D33").Formula = '="4.538"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  +2.61%  '

$ws.Range("D34").Formula = '="2.884"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  +3.06%  '

$ws.Range("D35").Formula = '="1.030"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  +2.08%  '

$ws.Range("E36").Value = '  +3.17%  '

$ws.Range("D37").Formula = '="0.01977"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  +2.90%  '

$ws.Range("D38").Formula = '="0.05319"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  +2.07%  '

$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Formula = '="0.5154"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +1.45%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Formula = '="2.814"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  +1.85%  '

$ws.Range("D41").Formula = '="0.1678"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  +2.31%  '

$ws.Range("D42").Formula = '="6.783"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +5.54%  '

$ws.Range("D43").Formula = '="8.505"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +2.92%  '

$ws.Range("D44").Formula = '="108.82"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +1.92%  '

$ws.Range("D45").Formula = '="10.59"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +1.89%  '

$ws.Range("D46").Formula = '="1.714"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +3.81%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Formula = '="0.4655"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  +2.14%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Formula = '="0.06407"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +1.99%  '

$ws.Range("D49").Formula = '="1.854"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  +2.95%  '

$ws.Range("D50").Formula = '="39.35"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  +4.42%  '

$ws.Range("D51").Formula = '="64.10"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +0.65%  '

$excel.CutCopyMode = 0
